# Edit slide 1 ("Subtitle 2" placeholder, Shapes.Item(2)):
#   - add two new paragraphs after "Research Supervisor - Senthil Kumar T":
#       "In support with : Sulakshan Vajipayajula"
#       "Architect CTO Office, IMB Security Bangalore - svajipay@in.ibm.com "
#     (the email text carries a mailto: hyperlink)

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$shp = $s.Shapes.Item(2)
$tf = $shp.TextFrame

# --- Paragraph: "In support with : " / "Sulakshan" / " " / "Vajipayajula" ---
$tr = $tf.TextRange
$lastPara = $tr.Paragraphs($tr.Paragraphs().Count)
$null = $lastPara.InsertAfter("`rIn support with : ")

$tr = $tf.TextRange
$lastPara = $tr.Paragraphs($tr.Paragraphs().Count)
$null = $lastPara.InsertAfter("Sulakshan")

$tr = $tf.TextRange
$lastPara = $tr.Paragraphs($tr.Paragraphs().Count)
$null = $lastPara.InsertAfter(" ")

$tr = $tf.TextRange
$lastPara = $tr.Paragraphs($tr.Paragraphs().Count)
$null = $lastPara.InsertAfter("Vajipayajula")

# --- Paragraph: "Architect CTO Office, IMB Security Bangalore <en-dash> " / "svajipay@in.ibm.com" / " " ---
$tr = $tf.TextRange
$lastPara = $tr.Paragraphs($tr.Paragraphs().Count)
$null = $lastPara.InsertAfter("`rArchitect CTO Office, IMB Security Bangalore $([char]0x2013) ")

$tr = $tf.TextRange
$lastPara = $tr.Paragraphs($tr.Paragraphs().Count)
$prefixStart = $lastPara.Start
$prefixLength = $lastPara.Length
$null = $lastPara.InsertAfter("svajipay@in.ibm.com")

$tr = $tf.TextRange
$lastPara = $tr.Paragraphs($tr.Paragraphs().Count)
$null = $lastPara.InsertAfter(" ")

# Select exactly the e-mail run (by character position) and turn it into a hyperlink.
$emailAddress = "svajipay@in.ibm.com"
$emailStart = $prefixStart + $prefixLength
$tr = $tf.TextRange
$emailRange = $tr.Characters($emailStart, $emailAddress.Length)
$hyperlink = $emailRange.ActionSettings(1).Hyperlink
$hyperlink.Address = "mailto:" + $emailAddress
